$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 103, shifting existing rows 103:134 down to 104:135.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with this week's data point.
$ws.Range("A103").Value = 11
$ws.Range("B103").Value = "Vega Monumental Concepción"
$ws.Range("C103").Value = "Bíobío"
$ws.Range("D103").Value = 44588
$ws.Range("E103").Value = 8
$ws.Range("F103").Value = 100112003
$ws.Range("G103").Value = "Ajo"
$ws.Range("H103").Value = "Chino"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 270
$ws.Range("K103").Value = 16000
$ws.Range("L103").Value = 17000
$ws.Range("M103").Value = 16444
$ws.Range("N103").Value = "`$/caja 10 kilos"
$ws.Range("O103").Value = "China"
$ws.Range("P103").Value = 1644
$ws.Range("Q103").Value = 10
$ws.Range("R103").Value = "Hortaliza"
